# Merge de cambios - actualizacion de metricas (fila 9: nuevos datos de iteracion
# "Metodos Norma 1, 2, inf en vector"), ancho de columna G y celda seleccionada.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Datos corregidos de la fila 9 (LoC real y tiempo de fin) ---
$ws.Range("C9").Value = 24
$ws.Range("F9").Value = 0.73611111111111116

# --- Ancho manual de la columna G (ya no es "best fit") ---
# El motor cuantiza ColumnWidth a pasos de 1/6 de caracter; 10.67 es el valor
# de entrada mas cercano al ancho final objetivo (11.42578125).
$ws.Columns("G").ColumnWidth = 10.67

# --- Celda seleccionada en la vista de la hoja ---
$ws.Range("H9").Select()

Write-Host "Metricas actualizadas"
